$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the two new use-case rows beneath the existing data
$ws.Range("B9").Value = "Gordon Ramsey best restaurants"
$ws.Range("B10").Value = "Couples date locations and fav romatic spots"

# Update the active selection to reflect the new last-row position
$ws.Range("B11").Select()
